{"js": "// Update CV content: rewrite the Professional Summary, rework the\n// Standard Chartered Bank bullet points (incl. a new bullet about the\n// credit-risk models), and rework the Think Big Analytics bullet points\n// (incl. a new blank bullet and a new \"Managed Hadoop clusters\" bullet\n// replacing the old fraud-detection one).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// ---- Helper: find the single paragraph whose text starts with a marker ----\nfunction findParagraph(marker) {\n  for (const p of paragraphs.items) {\n    if (p.text.indexOf(marker) === 0) {\n      return p;\n    }\n  }\n  throw new Error(\"Paragraph not found for marker: \" + marker);\n}\n\n// 1) Professional Summary paragraph.\nconst summaryPara = findParagraph(\n  \"Senior Engineering Leader with 15+ years bridging AI research\"\n);\n\n// 2) Standard Chartered Bank bullets.\nconst scLedPara = findParagraph(\n  \"Led enterprise AI and data platform development serving 11 markets\"\n);\nconst scDeliveredPara = findParagraph(\n  \"Delivered Self-Service ML Platform reducing model development\"\n);\nconst scModernizedPara = findParagraph(\n  \"Modernized MarTech infrastructure driving 30% increase\"\n);\n\n// 3) Think Big Analytics bullets.\nconst tbaEngineeredPara = findParagraph(\n  \"Engineered 5 data lakes processing 1.2 PB/hour with 20% optimization\"\n);\nconst tbaFraudPara = findParagraph(\n  \"Built real-time fraud detection reducing false positives by 60%\"\n);\n\n// ---- Apply the text edits --------------------------------------------------\n\nsummaryPara.insertText(\n  \"Senior Engineering Leader with 15+ years building teams, frameworks, and systems that turn Data and AI from research to production. Head of Data & Analytics for Google Cloud in Southeast Asia - a practice built from zero, delivering enterprise transformation across 7 countries. \\\"Player-Coach\\\" leading petabyte-scale platforms while driving published research (5 technical disclosures, 6 packages on PyPI/Maven, open-source AI safety tools).\",\n  Word.InsertLocation.replace\n);\n\nscLedPara.insertText(\n  \"Led enterprise AI and data platform transformation for retail banking.\",\n  Word.InsertLocation.replace\n);\n\nscDeliveredPara.insertText(\n  \"Built data & analytics platform serving 11 markets, 100+ systems, and 1200+ users\",\n  Word.InsertLocation.replace\n);\n\nscModernizedPara.insertText(\n  \"Delivered Self-Service ML Workbench reducing model deployment from months to weeks\",\n  Word.InsertLocation.replace\n);\n\n// New 4th Standard Chartered Bank bullet, right after the bullet above\n// (same ListParagraph/bullet formatting is inherited from scModernizedPara).\nscModernizedPara.insertParagraph(\n  \"Built credit risk models over 15K+ entities using news/social signals, reducing losses by $5M\",\n  Word.InsertLocation.after\n);\n\ntbaEngineeredPara.insertText(\n  \"Designed 5 data lakes processing 1.2 PB/hour and 40K daily files\",\n  Word.InsertLocation.replace\n);\n\n// New blank bullet right after the \"Designed 5 data lakes...\" bullet.\ntbaEngineeredPara.insertParagraph(\"\", Word.InsertLocation.after);\n\n// Replace the former \"Built real-time fraud detection...\" bullet text.\ntbaFraudPara.insertText(\n  \"Managed Hadoop clusters (300+ nodes) for banks and telcos across JAPAC\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# Update CV content: rewrite the Professional Summary, rework the\n# Standard Chartered Bank bullet points (incl. a new bullet about the\n# credit-risk models), and rework the Think Big Analytics bullet points\n# (incl. a new blank bullet and a new \"Managed Hadoop clusters\" bullet\n# replacing the old fraud-detection one).\n#\n# NOTE: paragraph *object references* can go stale once the paragraph\n# collection is mutated (e.g. after InsertParagraphAfter), so every\n# operation below re-resolves the target paragraph by its 1-based index\n# via $d.Paragraphs.Item($i) right before using it.\n\nfunction Find-ParagraphIndexByPrefix($doc, $prefix) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        if ($doc.Paragraphs.Item($i).Range.Text.StartsWith($prefix)) {\n            return $i\n        }\n    }\n    return -1\n}\n\n$d = $word.ActiveDocument\n\n# ---- 1) Professional Summary paragraph -------------------------------------\n$idx = Find-ParagraphIndexByPrefix $d \"Senior Engineering Leader with 15+ years bridging AI research\"\n$d.Paragraphs.Item($idx).Range.Text = 'Senior Engineering Leader with 15+ years building teams, frameworks, and systems that turn Data and AI from research to production. Head of Data & Analytics for Google Cloud in Southeast Asia - a practice built from zero, delivering enterprise transformation across 7 countries. \"Player-Coach\" leading petabyte-scale platforms while driving published research (5 technical disclosures, 6 packages on PyPI/Maven, open-source AI safety tools).'\n\n# ---- 2) Standard Chartered Bank bullets (plain text swaps) -----------------\n$idx = Find-ParagraphIndexByPrefix $d \"Led enterprise AI and data platform development serving 11 markets\"\n$d.Paragraphs.Item($idx).Range.Text = \"Led enterprise AI and data platform transformation for retail banking.\"\n\n$idx = Find-ParagraphIndexByPrefix $d \"Delivered Self-Service ML Platform reducing model development\"\n$d.Paragraphs.Item($idx).Range.Text = \"Built data & analytics platform serving 11 markets, 100+ systems, and 1200+ users\"\n\n$idx = Find-ParagraphIndexByPrefix $d \"Modernized MarTech infrastructure driving 30% increase\"\n$d.Paragraphs.Item($idx).Range.Text = \"Delivered Self-Service ML Workbench reducing model deployment from months to weeks\"\n\n# New 4th Standard Chartered Bank bullet, right after the bullet above.\n# InsertParagraphAfter() clones the paragraph's formatting (ListParagraph\n# style, numPr bullet, spacing), so the new paragraph is already a bullet.\n$idxDelivered = Find-ParagraphIndexByPrefix $d \"Delivered Self-Service ML Workbench reducing model deployment\"\n$d.Paragraphs.Item($idxDelivered).Range.InsertParagraphAfter()\n$idxNewCredit = $idxDelivered + 1\n$d.Paragraphs.Item($idxNewCredit).Range.Text = \"Built credit risk models over 15K+ entities using news/social signals, reducing losses by `$5M\"\n\n# ---- 3) Think Big Analytics bullets -----------------------------------------\n$idx = Find-ParagraphIndexByPrefix $d \"Engineered 5 data lakes processing 1.2 PB/hour with 20% optimization\"\n$d.Paragraphs.Item($idx).Range.Text = \"Designed 5 data lakes processing 1.2 PB/hour and 40K daily files\"\n\n# New blank bullet right after the \"Designed 5 data lakes...\" bullet.\n$idxDesigned = Find-ParagraphIndexByPrefix $d \"Designed 5 data lakes processing 1.2 PB/hour and 40K daily files\"\n$d.Paragraphs.Item($idxDesigned).Range.InsertParagraphAfter()\n\n# Replace the former \"Built real-time fraud detection...\" bullet text\n# (re-resolve by prefix since the new blank bullet shifted its index).\n$idxFraud = Find-ParagraphIndexByPrefix $d \"Built real-time fraud detection reducing false positives by 60%\"\n$d.Paragraphs.Item($idxFraud).Range.Text = \"Managed Hadoop clusters (300+ nodes) for banks and telcos across JAPAC\"\n"}
